{"js": "// The author was editing the title \"template\" (placing/leaving the cursor\n// after \"templa\") and finished the sentence in the \"R Markdown\" section.\n// As a result Word's automatic \"_GoBack\" bookmark (which marks the location\n// of the last edit) moved from the FirstParagraph text to inside the\n// heading word, and the two runs that the bookmark used to separate in the\n// FirstParagraph text were merged back into a single run.\n\n// --- 1) Remove the old \"_GoBack\" bookmark in the FirstParagraph text ---\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- 2) Merge the two runs it used to separate into a single run ---\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nconst firstParagraph = paragraphs.items[4]; // \"FirstParagraph\" styled paragraph\n\nconst beforeResults = firstParagraph.search(\n  \"This is an R Markdown document. Markdown is a simple formatting syntax for authoring \",\n  { matchCase: true }\n);\nbeforeResults.load(\"items\");\nconst afterResults = firstParagraph.search(\n  \"HTML, PDF, and MS Word documents. For more details on using R Markdown see \",\n  { matchCase: true }\n);\nafterResults.load(\"items\");\nawait context.sync();\n\nconst mergedRange = beforeResults.items[0].expandTo(afterResults.items[0]);\nmergedRange.insertText(\n  \"This is an R Markdown document. Markdown is a simple formatting syntax for authoring HTML, PDF, and MS Word documents. For more details on using R Markdown see \",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// --- 3) Split \"template\" into \"templa\" + \"te\" with a new _GoBack bookmark ---\nconst paragraphs2 = context.document.body.paragraphs;\nparagraphs2.load(\"items\");\nawait context.sync();\nconst heading = paragraphs2.items[0]; // \"Heading4\" styled paragraph (\"template\")\n\nconst templaResults = heading.search(\"templa\", { matchCase: true, matchWholeWord: false });\ntemplaResults.load(\"items\");\nawait context.sync();\n\nconst templaRange = templaResults.items[0];\nconst afterTemplaRange = templaRange.getRange(\"After\");\nafterTemplaRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# The author was finishing up the \"R Markdown\" section and left the\n# cursor inside the title word \"template\" (between \"templa\" and \"te\").\n# Word tracks the location of the user's last edit with a hidden\n# \"_GoBack\" bookmark; since that bookmark name is unique, re-placing it\n# both removes it from its old spot (in the FirstParagraph text, merging\n# the two runs it used to separate back into one run) and creates it at\n# the new spot (inside the heading), splitting that run in two.\n\n$d = $word.ActiveDocument\n\n# --- 1) Merge the FirstParagraph text back into a single run, removing\n#        the \"_GoBack\" bookmark that used to sit in the middle of it. ---\n$firstParagraphText = \"This is an R Markdown document. Markdown is a simple formatting syntax for authoring HTML, PDF, and MS Word documents. For more details on using R Markdown see \"\n$firstParagraphRange = $d.Paragraphs(5).Range\n$firstParagraphRange.Find.Execute(\n    $firstParagraphText, $false, $false, $false, $false, $false,\n    $true, 1, $false, $firstParagraphText, 2)\n\n# --- 2) Split \"template\" into \"templa\" + \"te\", placing a new \"_GoBack\"\n#        bookmark right after \"templa\". ---\n$headingRange = $d.Content\n$headingRange.Find.Execute(\"templa\")\n$headingRange.Collapse(0) # wdCollapseEnd\n$headingRange.Bookmarks.Add(\"_GoBack\")\n"}
